$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: name, email, repo link
$ws.Range("A3").Value = "حبيبة رزق عدلي "
$ws.Range("B3").Value = "habibarezk54@gmail.com"
$ws.Range("C3").Value = "https://github.com/WalTeR-RE/Open-Source-Uni-Project.git"

# Only the email (B3) gets a clickable mailto hyperlink + the same
# "Hyperlink" look the existing B2/C2 cells have.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:habibarezk54@gmail.com")
$ws.Range("B3").Style = $ws.Range("B2").Style

# Move the active selection to C3 (matches the saved workbook view).
$null = $ws.Range("C3").Select()
